$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder first two weeks of material:
# Row 3 (week for cm002) becomes the "grammar of graphics" topic (renamed)
# Row 5 (week for cm004) becomes "Data objects/tidy data"
$ws.Range("D3").Value = "Visualizations and the grammar of graphics"
$ws.Range("D5").Value = "Data objects/tidy data"

# Mark week 1 (row 2) as linked
$ws.Range("C2").Value = $true

# Update the active selection to D4
$ws.Range("D4").Select()
